$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "Federico Speroni"

$ws.Range("B43").Copy()
$ws.Range("B44").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B44").Value = 42871

$ws.Range("C44").Value = 5
$ws.Range("D44").Value = "Sprint 3 - BackEnd"
$ws.Range("E44").Value = "Alta Cliente, Actualizar Cliente, Actualizar Contrase;a, Obtener todos los clientes"

$ws.Range("E44").Select()
